$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New leading columns: last name, first name, position (plain text, no
# risk of being auto-detected as a number/date).
$ws.Range("A1").Value = "Rawls"
$ws.Range("B1").Value = "Thomas"
$ws.Range("C1").Value = "RB"

# Former columns B:I shift right by three (B->D, C->E, D->F, E->G, F->H,
# G->I, H->J, I->K). Values that look numeric/date-like need a leading
# quote so Excel stores them as literal text instead of converting them;
# the style is reset right after so no stray formatting sticks around.
$ws.Range("D1").Value = "'2018-10-07"
$ws.Range("D1").Style = "Normal"

$ws.Range("E1").Value = "'5"
$ws.Range("E1").Style = "Normal"

$ws.Range("F1").Value = "'25.065"
$ws.Range("F1").Style = "Normal"

$ws.Range("G1").Value = "CIN"

$ws.Range("H1").Value = "'"
$ws.Range("H1").Style = "Normal"

$ws.Range("I1").Value = "MIA"
$ws.Range("J1").Value = "W 27-17"

$ws.Range("K1").Value = "'"
$ws.Range("K1").Style = "Normal"

# New trailing numeric column.
$ws.Range("L1").Value = 0
